# Rename the header "idCommune" column to "stationId".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "stationId"

# The station-id column (A2:A34) was numeric ("number stored"); force it to
# be text (matches the xlsx's t="str" cell type for these rows) while
# keeping the same visible value "79049004".
$ws.Range("A2:A34").NumberFormat = "@"
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 1).Value = "79049004"
}

# Temperature columns D (temperatureMin) and E (temperatureMax) switch their
# decimal separator from "." to "," (e.g. "-4.7" -> "-4,7"). Values that have
# no decimal point (e.g. "9", "13") are left untouched.
for ($r = 2; $r -le 34; $r++) {
    foreach ($col in 4, 5) {
        $cell = $ws.Cells.Item($r, $col)
        $text = $cell.Text
        if ($text -and $text.Contains(".")) {
            $cell.Value = $text.Replace(".", ",")
        }
    }
}
